# Vantage Health Check On a Page -- "Update to VHC 4.1"
#
# Applies the visible-content edits from the authored diff:
#   - Title placeholder: "val:oap" -> "val:vhc" (both template tokens)
#   - "CPU Consumption" label -> "Busiest Hours of Week (CPU):"
#   - "COD:" label -> "Active CPU:" (and its auto-fit textbox widens)
#
# (Internal bookkeeping noise in the diff -- the slide's p:sldId number,
# the p14:creationId GUID-ish value, xmlns attribute ordering inside an
# a14:hiddenLine extension, the customXml part renumbering, and the
# cached text of the datetimeFigureOut footer fields -- are not exposed
# as writable PowerPoint Object Model surface, so they're intentionally
# left alone rather than hand-crafting XML outside the COM API.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) Title: "VHC on a page -- {{val:oap--intro.csv[1:2]}}: {{val:oap--intro.csv[1:1]}}"
#    becomes
#    "VHC on a page -- {{val:vhc--intro.csv[1:2]}}: {{val:vhc--intro.csv[1:1]}}"
# ---------------------------------------------------------------------
$title = Get-ShapeById $s 4
$titleRange = $title.TextFrame.TextRange

$full = $titleRange.Text
$firstOap = $full.IndexOf("val:oap")
$firstSub = $titleRange.Characters($firstOap + 1, 7)
$firstSub.Text = "val:vhc"

# re-typing "--intro.csv[1:2]}}: {{" as one contiguous edit merges the
# three runs it used to be split across into a single run
$full = $titleRange.Text
$markerIdx = $full.IndexOf("--intro.csv[1:2]}}: {{")
$markerLen = "--intro.csv[1:2]}}: {{".Length
$markerSub = $titleRange.Characters($markerIdx + 1, $markerLen)
$markerSub.Text = "--intro.csv[1:2]}}: {{"

$full = $titleRange.Text
$secondOap = $full.IndexOf("val:oap", $firstOap + 1)
$valSub = $titleRange.Characters($secondOap + 1, 4)
$valSub.Text = "val:"
$vhcSub = $titleRange.Characters($secondOap + 5, 3)
$vhcSub.Text = "vhc"

# ---------------------------------------------------------------------
# 2) "CPU Consumption" -> "Busiest Hours of Week (CPU):"
# ---------------------------------------------------------------------
$cpuLabel = Get-ShapeById $s 352
$cpuLabel.TextFrame.TextRange.Text = "Busiest Hours of Week (CPU):"

# ---------------------------------------------------------------------
# 3) "COD:" -> "Active CPU:" (textbox is wrap="none" + auto-fit, so it
#    grows to fit the new, longer caption)
# ---------------------------------------------------------------------
$codLabel = Get-ShapeById $s 79
$codLabel.TextFrame.TextRange.Text = "Active CPU:"
$codLabel.Width = 69.169
